$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - 18 updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3652
$ws1.Range("F5").Value = 3652
$ws1.Range("F7").Value = 5179
$ws1.Range("F9").Value = 381
$ws1.Range("F16").Value = 325
$ws1.Range("F17").Value = 39
$ws1.Range("F21").Value = 364
$ws1.Range("F22").Value = 4954
$ws1.Range("F23").Value = 46
$ws1.Range("F26").Value = 6081
$ws1.Range("F28").Value = 18
$ws1.Range("F29").Value = 3235
$ws1.Range("F30").Value = 349
$ws1.Range("F34").Value = 128
$ws1.Range("F36").Value = 1064
$ws1.Range("F37").Value = 86
$ws1.Range("F40").Value = 887
$ws1.Range("F41").Value = 1043

# Sheet "演出" (sheet2) - 1 update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 61

# Sheet "本地生活" (sheet3) - 1 update
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1130

# Sheet "全部类型" (sheet4) - 20 updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1130
$ws4.Range("F7").Value = 3653
$ws4.Range("F8").Value = 3653
$ws4.Range("F10").Value = 5179
$ws4.Range("F12").Value = 381
$ws4.Range("F19").Value = 325
$ws4.Range("F20").Value = 39
$ws4.Range("F25").Value = 364
$ws4.Range("F26").Value = 4954
$ws4.Range("F27").Value = 46
$ws4.Range("F30").Value = 6081
$ws4.Range("F32").Value = 18
$ws4.Range("F33").Value = 3235
$ws4.Range("F34").Value = 349
$ws4.Range("F39").Value = 128
$ws4.Range("F41").Value = 1064
$ws4.Range("F42").Value = 86
$ws4.Range("F45").Value = 887
$ws4.Range("F46").Value = 1043
$ws4.Range("F50").Value = 61
